$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'41.228.30"
$ws.Range("E2").Value = "  -3.58%  "
$ws.Range("D3").Value = "'2.463.00"
$ws.Range("E3").Value = "  -2.58%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'312.08"
$ws.Range("E5").Value = "  +0.72%  "
$ws.Range("D6").Value = "'93.96"
$ws.Range("E6").Value = "  -6.19%  "
$ws.Range("E7").Value = "  -2.96%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'0.497"
$ws.Range("E9").Value = "  -4.77%  "
$ws.Range("D10").Value = "'33.40"
$ws.Range("E10").Value = "  -5.87%  "
$ws.Range("D11").Value = "'0.0778"
$ws.Range("E11").Value = "  -3.46%  "
$ws.Range("E12").Value = "  -1.22%  "
$ws.Range("D13").Value = "'6.99"
$ws.Range("E13").Value = "  -4.46%  "
$ws.Range("D14").Value = "'2.842.87"
$ws.Range("E14").Value = "  -2.56%  "
$ws.Range("D15").Value = "'2.453.36"
$ws.Range("E15").Value = "  -2.03%  "
$ws.Range("D16").Value = "'14.87"
$ws.Range("E16").Value = "  -3.02%  "
$ws.Range("D17").Value = "'0.782"
$ws.Range("E17").Value = "  -3.82%  "
$ws.Range("D18").Value = "'41.180.63"
$ws.Range("E19").Value = "  -5.54%  "
$ws.Range("E20").Value = "  -3.15%  "
$ws.Range("E21").Value = "  -8.26%  "
$ws.Range("D22").Value = "'68.32"
$ws.Range("E22").Value = "  -1.46%  "
$ws.Range("D23").Value = "'235.81"
$ws.Range("E23").Value = "  -2.99%  "
$ws.Range("E24").Value = "  -4.07%  "
$ws.Range("E25").Value = "  +0.22%  "
$ws.Range("D27").Value = "'24.02"
$ws.Range("E27").Value = "  -5.29%  "
$ws.Range("E28").Value = "  -5.67%  "
$ws.Range("D29").Value = "'9.62"
$ws.Range("E29").Value = "  -5.36%  "
$ws.Range("D30").Value = "'36.44"
$ws.Range("E30").Value = "  -5.37%  "
$ws.Range("D31").Value = "'151.88"
$ws.Range("E31").Value = "  -5.35%  "
$ws.Range("E32").Value = "  -4.91%  "
$ws.Range("E33").Value = "  -4.83%  "
$ws.Range("E34").Value = "  -4.30%  "
$ws.Range("D35").Value = "'0.0743"
$ws.Range("E35").Value = "  -5.21%  "
$ws.Range("D36").Value = "'3.05"
$ws.Range("E36").Value = "  -2.30%  "
$ws.Range("B37").Value = "Celestia"
$ws.Range("C37").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D37").Value = "'17.09"
$ws.Range("E37").Value = "  -7.28%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").Value = "'1.89"
$ws.Range("E38").Value = "  -3.95%  "
$ws.Range("E39").Value = "  -2.96%  "
$ws.Range("E40").Value = "  -7.78%  "
$ws.Range("D41").Value = "'4.26"
$ws.Range("E41").Value = "  +2.08%  "
$ws.Range("E42").Value = "  +0.17%  "
$ws.Range("D43").Value = "'20.05"
$ws.Range("E43").Value = "  -10.85%  "
$ws.Range("D44").Value = "'1.977.49"
$ws.Range("E44").Value = "  -1.48%  "
$ws.Range("E45").Value = "  -5.66%  "
$ws.Range("E46").Value = "  -8.02%  "
$ws.Range("D47").Value = "'8.69"
$ws.Range("E47").Value = "  -2.18%  "
$ws.Range("B48").Value = "ordi"
$ws.Range("C48").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D48").Value = "'69.09"
$ws.Range("E48").Value = "  -3.86%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "'96.95"
$ws.Range("E49").Value = "  -3.72%  "
$ws.Range("E50").Value = "  -6.60%  "
$ws.Range("D51").Value = "'74.17"
$ws.Range("E51").Value = "  -6.35%  "
